# Applies the "Ajout du changement de signature" commit:
#  - B1: "SFC 2020/2020" -> "SFC 2019/2020"
#  - H3: "du 01/01/2020 au 05/05/2020" -> "du 16/09/2019 au 20/06/2020"
#  - I6: "2020 T1" -> "2019 T1"
#  - Q14: 105 -> 94

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "SFC 2019/2020"
$ws.Range("H3").Value = "du 16/09/2019 au 20/06/2020"
$ws.Range("I6").Value = "2019 T1"
$ws.Range("Q14").Value = 94
